$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: course changes from FIN 372 to EC 480, with an expanded prerequisite chain ---
$ws.Range("B4").Value = "EC 480"
$ws.Range("D4").Value = "EC 301.0"
$ws.Range("E4").Value = "Missing"

$ws.Range("F4").Value = "prerequisite"
$ws.Range("G4").Value = "EC 302.0"
$ws.Range("H4").Value = "Grade"

$ws.Range("I4").Value = "prerequisite"
$ws.Range("J4").Value = "FIN 301.0"
$ws.Range("K4").Value = "Grade"

$ws.Range("L4").Value = "prerequisite"
$ws.Range("M4").Value = "MA 209.0"
$ws.Range("N4").Value = "Grade"

# New cells on row 4 need the same bordered style as the rest of the row
$ws.Range("F4:N4").Borders.LineStyle = 1

# --- Row 5: additional missing prerequisite line for EC 480 ---
$ws.Range("A5").Value = "Elettra Scianetti"
$ws.Range("B5").Value = "EC 480"
$ws.Range("L5").Value = "prerequisite"
$ws.Range("M5").Value = "EC 360.0"
$ws.Range("N5").Value = "Grade"

$ws.Range("A5:B5").Borders.LineStyle = 1
$ws.Range("L5:N5").Borders.LineStyle = 1

# --- Row 6: the FIN 372 / FIN 301.0 entry that used to live on row 4 ---
$ws.Range("A6").Value = "Elettra Scianetti"
$ws.Range("B6").Value = "FIN 372"
$ws.Range("C6").Value = "prerequisite"
$ws.Range("D6").Value = "FIN 301.0"
$ws.Range("E6").Value = "Grade"

$ws.Range("A6:E6").Borders.LineStyle = 1
